$d = $word.ActiveDocument

$pairs = @(
    @("81×58=4698", "84×12=1008"),
    @("45×71=3195", "14×20=280"),
    @("39×47=1833", "57×56=3192"),
    @("17×33=561", "89×28=2492"),
    @("35×13=455", "43×43=1849"),
    @("65×89=5785", "64×88=5632"),
    @("44×28=1232", "53×50=2650"),
    @("14×41=574", "70×35=2450"),
    @("47×29=1363", "36×77=2772"),
    @("16×26=416", "38×37=1406"),
    @("60×28=1680", "52×97=5044"),
    @("52×27=1404", "35×65=2275"),
    @("74×98=7252", "43×70=3010"),
    @("17×50=850", "58×99=5742"),
    @("52×19=988", "35×50=1750"),
    @("95×94=8930", "40×41=1640"),
    @("92×17=1564", "60×68=4080"),
    @("51×27=1377", "46×77=3542"),
    @("20×48=960", "49×23=1127"),
    @("18×86=1548", "79×11=869"),
    @("80×26=2080", "74×36=2664"),
    @("63×31=1953", "45×38=1710"),
    @("82×86=7052", "48×17=816"),
    @("12×88=1056", "20×18=360"),
    @("54×11=594", "42×13=546")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
